# Fixed sb_mean and sb_diff
# Column H holds SB_BinaryStats_mean_longstretch1 values. A batch of rows
# had been off-by-one (and flagged with the red "outlier" highlight style).
# Correct the values (decrement by 1) and clear the highlight back to the
# workbook's default (unstyled) cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4, 13, 15, 16, 18, 19, 20, 23, 24, 25, 30, 31, 35, 37, 40, 44, 49, 52, 54, 58, 60, 67)

foreach ($r in $rows) {
    $cell = $ws.Range("H$r")
    $current = $cell.Value2
    $cell.Value = $current - 1
    $cell.Style = "Normal"
}
